$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format D2:E51 as Text first so numeric-looking strings (prices, thousand-dot
# separated values, percentages) are preserved verbatim instead of Excel
# auto-converting them to numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "28.443.08"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.827.62"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "317.21"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").Value = "0.5343"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D9").Value = "0.07608"
$ws.Range("E9").Value = "  +2.24%  "
$ws.Range("D10").Value = "41.80"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("D11").Value = "1.110"
$ws.Range("E11").Value = "  +1.62%  "
$ws.Range("D12").Value = "6.329"
$ws.Range("E12").Value = "  +4.15%  "
$ws.Range("D13").Value = "7.619"
$ws.Range("E13").Value = "  +5.55%  "
$ws.Range("D14").Value = "1.000"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "20.96"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").Value = "1.828.72"
$ws.Range("E16").Value = "  +3.21%  "
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "0.00001074"
$ws.Range("E18").Value = "  +2.17%  "
$ws.Range("D19").Value = "0.06586"
$ws.Range("E19").Value = "  +2.34%  "
$ws.Range("D21").Value = "0.9997"
$ws.Range("E21").Value = "  -0.01%  "
$ws.Range("D22").Value = "6.083"
$ws.Range("E22").Value = "  +3.62%  "
$ws.Range("D23").Value = "28.452.25"
$ws.Range("E23").Value = "  +1.53%  "
$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("E25").Value = "  +6.12%  "
$ws.Range("D26").Value = "2.456"
$ws.Range("E26").Value = "  +7.70%  "
$ws.Range("D27").Value = "157.28"
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "20.66"
$ws.Range("E28").Value = "  +2.24%  "
$ws.Range("D29").Value = "2.043.17"
$ws.Range("E29").Value = "  +3.29%  "
$ws.Range("D30").Value = "124.05"
$ws.Range("E30").Value = "  +3.43%  "
$ws.Range("D31").Value = "1.123"
$ws.Range("E31").Value = "  +1.74%  "
$ws.Range("D32").Value = "0.1102"
$ws.Range("E32").Value = "  +4.59%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "5.666"
$ws.Range("E33").Value = "  +2.84%  "
$ws.Range("B34").Value = "Hedera"
$ws.Range("C34").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "0.07477"
$ws.Range("E34").Value = "  +16.20%  "
$ws.Range("D35").Value = "3.643"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "0.2227"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("E37").Value = "  +2.80%  "
$ws.Range("E38").Value = "  +4.26%  "
$ws.Range("D39").Value = "8.843"
$ws.Range("E39").Value = "  +5.13%  "
$ws.Range("D40").Value = "0.6260"
$ws.Range("E40").Value = "  +2.06%  "
$ws.Range("D41").Value = "11.30"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").Value = "1.176"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "0.9993"
$ws.Range("E43").Value = "  -0.02%  "
$ws.Range("E44").Value = "  -3.16%  "
$ws.Range("D45").Value = "13.53"
$ws.Range("E45").Value = "  +1.70%  "
$ws.Range("D46").Value = "3.702"
$ws.Range("E46").Value = "  +0.91%  "
$ws.Range("D47").Value = "0.5844"
$ws.Range("E47").Value = "  +1.86%  "
$ws.Range("D48").Value = "124.83"
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D49").Value = "2.001"
$ws.Range("E49").Value = "  +4.18%  "
$ws.Range("E50").Value = "  +1.28%  "
$ws.Range("D51").Value = "0.06889"
$ws.Range("E51").Value = "  +1.44%  "

# Restore the original (unformatted/General) appearance now that the text
# values are committed, so no stray text-format style lingers on the cells.
$ws.Range("D2:E51").ClearFormats()

